$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 87; existing rows 87-116 shift down to 88-117.
$ws.Rows.Item(87).Insert()

# Populate the newly inserted row 87 with the new record.
$r = 87
$ws.Cells.Item($r, 1).Value = 4
$ws.Cells.Item($r, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($r, 3).Value = "Los Lagos"
$ws.Cells.Item($r, 4).Value = 44663
$ws.Cells.Item($r, 5).Value = 10
$ws.Cells.Item($r, 6).Value = 100112052
$ws.Cells.Item($r, 7).Value = "Albahaca"
$ws.Cells.Item($r, 8).Value = "Sin especificar"
$ws.Cells.Item($r, 9).Value = "Primera"
$ws.Cells.Item($r, 10).Value = 90
$ws.Cells.Item($r, 11).Value = 6000
$ws.Cells.Item($r, 12).Value = 6000
$ws.Cells.Item($r, 13).Value = 6000
$ws.Cells.Item($r, 14).Value = "`$/docena de matas"
$ws.Cells.Item($r, 15).Value = "Región Metropolitana"
$ws.Cells.Item($r, 16).Value = 1000
$ws.Cells.Item($r, 17).Value = 6
$ws.Cells.Item($r, 18).Value = "Hortaliza"
